$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 621.95
$ws.Range("C2").Value = 612.3
$ws.Range("D2").Value = 617.9
$ws.Range("E2").Value = 617.05
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 613.1
$ws.Range("B3").Value = 3037.9
$ws.Range("C3").Value = 3002
$ws.Range("D3").Value = 3012
$ws.Range("E3").Value = 3015.35
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 3013.75
$ws.Range("B4").Value = 515.5
$ws.Range("C4").Value = 509.5
$ws.Range("D4").Value = 511.3
$ws.Range("E4").Value = 512
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 510.15
$ws.Range("B5").Value = 1874.7
$ws.Range("C5").Value = 1850
$ws.Range("D5").Value = 1865.1
$ws.Range("E5").Value = 1864.95
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 1868
$ws.Range("B6").Value = 7286.65
$ws.Range("C6").Value = 7193
$ws.Range("D6").Value = 7235
$ws.Range("E6").Value = 7244.9
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 7275.5
$ws.Range("B7").Value = 208.25
$ws.Range("C7").Value = 202.51
$ws.Range("D7").Value = 203.7
$ws.Range("E7").Value = 203.66
$ws.Range("F7").Value = 273
$ws.Range("G7").Value = 202.67
$ws.Range("B8").Value = 244.45
$ws.Range("C8").Value = 240.85
$ws.Range("D8").Value = 244.05
$ws.Range("E8").Value = 243.85
$ws.Range("F8").Value = 141
$ws.Range("G8").Value = 242.75
$ws.Range("B9").Value = 505.2
$ws.Range("C9").Value = 495.15
$ws.Range("D9").Value = 497.2
$ws.Range("E9").Value = 497
$ws.Range("F9").Value = 108
$ws.Range("G9").Value = 502.55
$ws.Range("B10").Value = 851.8
$ws.Range("C10").Value = 840.15
$ws.Range("D10").Value = 841.5
$ws.Range("E10").Value = 841.65
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 850.6
$ws.Range("B11").Value = 4854.55
$ws.Range("C11").Value = 4787.15
$ws.Range("D11").Value = 4800
$ws.Range("E11").Value = 4800.1
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 4848.45
$ws.Range("B12").Value = 190.68
$ws.Range("C12").Value = 187.61
$ws.Range("D12").Value = 190
$ws.Range("E12").Value = 189.65
$ws.Range("F12").Value = 61
$ws.Range("G12").Value = 188.28
$ws.Range("B13").Value = 1797.7
$ws.Range("C13").Value = 1767.3
$ws.Range("D13").Value = 1794.6
$ws.Range("E13").Value = 1790.55
$ws.Range("F13").Value = 28
$ws.Range("G13").Value = 1780.05
$ws.Range("B14").Value = 1649.5
$ws.Range("C14").Value = 1642
$ws.Range("D14").Value = 1644.4
$ws.Range("E14").Value = 1645.45
$ws.Range("F14").Value = 114
$ws.Range("G14").Value = 1648.25
$ws.Range("B15").Value = 672.4
$ws.Range("C15").Value = 667.25
$ws.Range("D15").Value = 669.3
$ws.Range("E15").Value = 669.95
$ws.Range("F15").Value = 37
$ws.Range("G15").Value = 669.05
$ws.Range("B16").Value = 1241.85
$ws.Range("C16").Value = 1231.3
$ws.Range("D16").Value = 1235.35
$ws.Range("E16").Value = 1235.95
$ws.Range("F16").Value = 91
$ws.Range("G16").Value = 1241
$ws.Range("B17").Value = 1440.55
$ws.Range("C17").Value = 1418.9
$ws.Range("D17").Value = 1426.2
$ws.Range("E17").Value = 1422.9
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 1435.85
$ws.Range("B18").Value = 1940
$ws.Range("C18").Value = 1909.6
$ws.Range("D18").Value = 1933
$ws.Range("E18").Value = 1933.15
$ws.Range("F18").Value = 40
$ws.Range("G18").Value = 1926.9
$ws.Range("B19").Value = 955.65
$ws.Range("C19").Value = 942.6
$ws.Range("D19").Value = 949.8
$ws.Range("E19").Value = 950.5
$ws.Range("F19").Value = 7
$ws.Range("G19").Value = 950
$ws.Range("B20").Value = 719.5
$ws.Range("C20").Value = 700.6
$ws.Range("D20").Value = 707.75
$ws.Range("E20").Value = 707.4
$ws.Range("F20").Value = 61
$ws.Range("G20").Value = 703.25
$ws.Range("B21").Value = 2746.55
$ws.Range("C21").Value = 2711.9
$ws.Range("D21").Value = 2721
$ws.Range("E21").Value = 2723.1
$ws.Range("F21").Value = 23
$ws.Range("G21").Value = 2745.9
$ws.Range("B22").Value = 332.3
$ws.Range("C22").Value = 326
$ws.Range("D22").Value = 328.5
$ws.Range("E22").Value = 329.1
$ws.Range("F22").Value = 27
$ws.Range("G22").Value = 328.45
$ws.Range("B23").Value = 408.4
$ws.Range("C23").Value = 401.85
$ws.Range("D23").Value = 403
$ws.Range("E23").Value = 403.25
$ws.Range("F23").Value = 142
$ws.Range("G23").Value = 405.05
$ws.Range("B24").Value = 3041
$ws.Range("C24").Value = 2974.1
$ws.Range("D24").Value = 2991
$ws.Range("E24").Value = 2985.95
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 3039.65
$ws.Range("B25").Value = 822.15
$ws.Range("C25").Value = 815.6
$ws.Range("D25").Value = 818
$ws.Range("E25").Value = 818.75
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 816.35
$ws.Range("B26").Value = 812.5
$ws.Range("C26").Value = 794.9
$ws.Range("D26").Value = 798.5
$ws.Range("E26").Value = 796.4
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 807.1
$ws.Range("B27").Value = 1088.25
$ws.Range("C27").Value = 1076.3
$ws.Range("D27").Value = 1081.45
$ws.Range("E27").Value = 1081.25
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 1079.3
$ws.Range("B28").Value = 1084
$ws.Range("C28").Value = 1066.95
$ws.Range("D28").Value = 1068.7
$ws.Range("E28").Value = 1069.15
$ws.Range("F28").Value = 79
$ws.Range("G28").Value = 1080
$ws.Range("B29").Value = 424
$ws.Range("C29").Value = 419.15
$ws.Range("D29").Value = 420.7
$ws.Range("E29").Value = 420.9
$ws.Range("F29").Value = 81
$ws.Range("G29").Value = 421.3
$ws.Range("B30").Value = 152.79
$ws.Range("C30").Value = 151.4
$ws.Range("D30").Value = 152.03
$ws.Range("E30").Value = 151.72
$ws.Range("F30").Value = 283
$ws.Range("G30").Value = 152.66
$ws.Range("B31").Value = 11822.75
$ws.Range("C31").Value = 11514.2
$ws.Range("D31").Value = 11532
$ws.Range("E31").Value = 11542.65
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 11781.3
